{"js": "// The author expanded the description of the second tested class\n// (getRMAddress -> ClientRMProxy) and named the JUnit test file that was\n// produced for it, while leaving the rest of the sentence intact.\n//\n// Old: \"...named as AHSproxy and getRMAddress and done the Unit testing for\n//       these two classes. The logs and classes.txt are placed in the\n//       respective folder of Cyclomatic complexity.\"\n// New: \"...named as AHSproxy and ClientRMProxy done the Unit testing for\n//       these two classes. The logs and\n//       TestClientRMProxyCyclomaticComplexity.java(test class) are placed in\n//       the respective folder of Cyclomatic complexity.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that carries the sentence we need to edit.\nconst target = paragraphs.items.find(\n  (p) => p.text.indexOf(\"We have taken two classes from Release 3.0.0\") !== -1\n);\nif (!target) {\n  throw new Error(\"Could not find the target paragraph\");\n}\n\n// Replace the two phrases that changed; everything else in the sentence\n// (including the trailing period, which lives in its own run after a\n// bookmark) is left untouched.\nconst firstHit = target.search(\"getRMAddress and\", { matchCase: true });\nfirstHit.load(\"items\");\nawait context.sync();\nif (firstHit.items.length === 0) {\n  throw new Error(\"Could not find 'getRMAddress and' to replace\");\n}\nfirstHit.items[0].insertText(\"ClientRMProxy\", \"Replace\");\nawait context.sync();\n\nconst secondHit = target.search(\n  \"classes. The logs and classes.txt are placed\",\n  { matchCase: true }\n);\nsecondHit.load(\"items\");\nawait context.sync();\nif (secondHit.items.length === 0) {\n  throw new Error(\"Could not find 'classes. The logs and classes.txt are placed' to replace\");\n}\nsecondHit.items[0].insertText(\n  \"classes. The logs and TestClientRMProxyCyclomaticComplexity.java(test class) are placed\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# The author expanded the description of the second tested class\n# (getRMAddress -> ClientRMProxy) and named the JUnit test file that was\n# produced for it, while leaving the rest of the sentence intact.\n#\n# Old: \"...named as AHSproxy and getRMAddress and done the Unit testing for\n#       these two classes. The logs and classes.txt are placed in the\n#       respective folder of Cyclomatic complexity.\"\n# New: \"...named as AHSproxy and ClientRMProxy done the Unit testing for\n#       these two classes. The logs and\n#       TestClientRMProxyCyclomaticComplexity.java(test class) are placed in\n#       the respective folder of Cyclomatic complexity.\"\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceOne = 1\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Execute(\"getRMAddress and\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"ClientRMProxy\", $wdReplaceOne)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Execute(\"classes. The logs and classes.txt are placed\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"classes. The logs and TestClientRMProxyCyclomaticComplexity.java(test class) are placed\", $wdReplaceOne)\n"}
